# The "domains" sheet is stored internally as "site.conf" and is the
# sheet that was active/selected when the file was last edited.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("site.conf")
$ws.Activate()

# Column L held a "firmware/BSSID" style helper formula that appended a
# trailing slash to the domain code in column A (e.g. "ffhpd01/"). That
# trailing slash was never actually used downstream, so drop it: the
# formula becomes a straight reference to column A.
for ($row = 2; $row -le 25; $row++) {
    $ws.Cells.Item($row, 12).Formula = "=A$row"
}

# Reflect the reviewer's new selection/scroll position: column L
# (the parameter that was just edited) is selected, scrolled into view,
# with L2 left as the active cell inside that selection.
$ws.Range("L2:L25").Select()
$excel.ActiveWindow.ScrollColumn = 8
$excel.ActiveWindow.ScrollRow = 1

$wb.Save()
